# Commit: "changed names of transitions in Example 2 dataset changed respective
# names in template.xlsx added an additional pane for upload of different
# template.xlsx file"
#
# Net effect on this workbook: every transition name in the "Example2" sheet's
# header row (currently suffixed "_quan") gets renamed to end in "_quant", and
# the sheet's active-cell selection moves from G17 to J1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example2")

# Rename every transition column header ("..._quan" -> "..._quant") across
# the whole sheet via a literal find & replace (mirrors renaming the
# corresponding entries in templates.xlsx).
$ws.Cells.Replace("_quan", "_quant")

# Move the active selection to J1 (previously G17).
$ws.Activate()
$ws.Range("J1").Select()
